$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.041.44'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.793.07'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.16%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.93'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.79%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.08%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5354'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -1.83%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.54%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07425'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.98%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.91'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.73%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.091'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -2.90%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.00%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.55'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.77%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.120'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.02%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.239'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -2.22%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.787.34'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.45%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.70%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001057'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.36%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06499'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.84%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.00%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.25'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.31%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.888'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.15%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.066.32'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.20%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -2.60%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.86%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '155.44'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.61%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.94%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.992.05'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.52%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -3.91%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '121.17'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.73%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.78%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +3.19%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.664'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.42%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.556'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -3.44%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.2247'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -4.33%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06481'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -4.28%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02288'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.55%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -2.83%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.452'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -3.63%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6188'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -3.25%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.56%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.12'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -4.62%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.179'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +1.74%  '
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.35'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.26%  '
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = 'PancakeSwap'
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.669'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.11%  '
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5779'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -3.29%  '
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '124.81'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.34%  '
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'EOS'
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.188'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +3.23%  '
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.926'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -3.83%  '
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06822'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.92%  '
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '71.21'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -2.33%  '
